$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix up the intro paragraph text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "These are the results of the most previous doubles match with ID: 2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "These are the results of the most previous singles match with ID: 1", 2)

$d.Content.Find.Execute(
    "Member/Team 1 had a final score of - 3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Member/Team 1 had a final score of - 21", 2)

$d.Content.Find.Execute(
    "Member/Team 2 had a final score of - 21",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Member/Team 2 had a final score of - 8", 2)

# ---------------------------------------------------------------------------
# 2) Fix up / extend the results table.
#    The first data row (header) plus the next 7 data rows (rows 2-8 of the
#    table) are unchanged. From row 9 onward the score progression differs,
#    and 5 brand-new rows are appended at the end.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Grow the table to the final size (31 rows total = 1 header + 30 data rows).
while ($t.Rows.Count -lt 31) {
    $t.Rows.Add() | Out-Null
}

$scores = @(
    @(4,3),
    @(5,3),
    @(6,3),
    @(6,4),
    @(7,4),
    @(7,5),
    @(8,5),
    @(9,5),
    @(10,5),
    @(10,6),
    @(10,7),
    @(11,7),
    @(12,7),
    @(12,8),
    @(13,8),
    @(14,8),
    @(15,8),
    @(16,8),
    @(17,8),
    @(18,8),
    @(19,8),
    @(20,8),
    @(21,8)
)

$startRow = 9
for ($i = 0; $i -lt $scores.Length; $i++) {
    $rowIndex = $startRow + $i
    $pair = $scores[$i]
    $t.Cell($rowIndex, 1).Range.Text = [string]$pair[0]
    $t.Cell($rowIndex, 2).Range.Text = [string]$pair[1]
}
